$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up number formats on existing rows (7-9, 12-14) ---
# P7 was using the soon-to-be-removed 0.00000 custom format -> 0.00
$ws.Range("P7").NumberFormat = "0.00"
# P8 / P9 had no explicit format (General) -> 0.00
$ws.Range("P8").NumberFormat = "0.00"
$ws.Range("P9").NumberFormat = "0.00"
# N9 / O9 move from 0.00 to the integer "0" format
$ws.Range("N9").NumberFormat = "0"
$ws.Range("O9").NumberFormat = "0"
# P12:P14 move from the 0.000000 custom format to 0.00
$ws.Range("P12").NumberFormat = "0.00"
$ws.Range("P13").NumberFormat = "0.00"
$ws.Range("P14").NumberFormat = "0.00"

# --- Fill in new row 10 data (Baseline 2010 C81) ---
$ws.Range("A10").Value = "CW3M"
$ws.Range("B10").Value = "Baseline 2010 C81"
$ws.Range("C10").ClearFormats()
$ws.Range("C10").Value = 2010
$ws.Range("D10").Value = 1034.060303
$ws.Range("E10").Value = 1990.4676509999999
$ws.Range("F10").Value = 1.255063
$ws.Range("G10").Value = 327.58108499999997
$ws.Range("H10").Value = 10.610913999999999
$ws.Range("I10").Value = 8.8404570000000007
$ws.Range("J10").Value = 814.38360599999999
$ws.Range("K10").Value = 93.229797000000005
$ws.Range("L10").Value = 1292.8286129999999
$ws.Range("M10").Value = 1158.413818
$ws.Range("N10").NumberFormat = "0"
$ws.Range("N10").Value = 7105.0297849999997
$ws.Range("O10").NumberFormat = "0"
$ws.Range("O10").Value = 29450.638672000001
$ws.Range("P10").NumberFormat = "0.00"
$ws.Range("P10").Value = 3.7212749999999999
$ws.Range("Q10").ClearFormats()
$ws.Range("Q10").Value = 0.001106
$ws.Range("R10").Value = 2010

# --- Add new row 15 data (Baseline 2010-18 C81) ---
$ws.Range("A15").Value = "CW3M"
$ws.Range("B15").Value = "Baseline 2010-18 C81"
$ws.Range("C15").Value = "2010-18"
$ws.Range("D15:M15").NumberFormat = "0.00"
$ws.Range("D15").Value = 1135.8478461111113
$ws.Range("E15").Value = 1901.5157334444443
$ws.Range("F15").Value = 1.0119255555555557
$ws.Range("G15").Value = 327.78053433333326
$ws.Range("H15").Value = 9.775355222222224
$ws.Range("I15").Value = 8.145128999999999
$ws.Range("J15").Value = 769.26112866666654
$ws.Range("K15").Value = 83.47062044444445
$ws.Range("L15").Value = 1378.3211942222222
$ws.Range("M15").Value = 1141.5044894444445
$ws.Range("N15").NumberFormat = "0"
$ws.Range("N15").Value = 4878.4023980000002
$ws.Range("O15").NumberFormat = "0"
$ws.Range("O15").Value = 27227.338324888889
$ws.Range("P15").NumberFormat = "0.00"
$ws.Range("P15").Value = 4.7711666666666668
$ws.Range("Q15").NumberFormat = "0.000000"
$ws.Range("Q15").Value = 0.0014151111111111109
$ws.Range("R15").Value = "2010-18"

# --- Update the selection to match the saved view ---
$ws.Range("P12:P14").Select()
